$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.401.19"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "2.637.19"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'598.08"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'152.28"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  +5.39%  "
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("D11").Value = "'0.395"
$ws.Range("E11").Value = "  +3.45%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("D13").Value = "'28.12"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "3.102.09"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("E15").Value = "  +14.65%  "
$ws.Range("D16").Value = "64.192.30"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "2.634.05"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "'12.31"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'4.80"
$ws.Range("D20").Value = "'350.62"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'67.74"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").Value = "'9.24"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'8.32"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'556.71"
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "'0.162"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "0.0₃0915"
$ws.Range("E31").Value = "  +8.38%  "
$ws.Range("D32").Value = "'2.08"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("E33").Value = "  +5.43%  "
$ws.Range("D34").Value = "'5.50"
$ws.Range("E34").Value = "  +4.53%  "
$ws.Range("D35").Value = "'6.21"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +3.66%  "
$ws.Range("D37").Value = "'165.89"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").Value = "'20.15"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D42").Value = "'169.66"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("E43").Value = "  +4.85%  "
$ws.Range("D44").Value = "'23.31"
$ws.Range("E44").Value = "  +9.06%  "
$ws.Range("E45").Value = "  +11.99%  "
$ws.Range("D46").Value = "'0.0589"
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("D48").Value = "'0.0253"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "'19.39"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  +19.49%  "
